$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 8, pushing existing rows 8.. down by one.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row with the new "all" / "GBD region names" entry.
$ws.Range("A8").Value = "all"
$ws.Range("B8").Value = "GBD region names"

# Match formatting style of the other "all" rows (right-aligned column A).
$ws.Range("A8").HorizontalAlignment = -4152

# Update the active selection to reflect where the edit was made.
$ws.Range("A9").Select()
